# Apply "change predict income algorithm": the balance-sheet roll-forward now
# drops the oldest reported fiscal period (1396/12) and appends the newest one
# (1401/12), shifting every existing period one column to the left (D<-E, E<-F,
# F<-G, G<-H) and filling H with the newly predicted/reported period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: fiscal-period column headers ---
$headerRow8 = @("12 ماهه منتهی به 1397/12", "12 ماهه منتهی به 1398/12", "12 ماهه منتهی به 1399/12", "12 ماهه منتهی به 1400/12", "12 ماهه منتهی به 1401/12")
for ($i = 0; $i -lt $headerRow8.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $headerRow8[$i]
}

# --- Row 9: publish-date column headers ---
$headerRow9 = @("1399-04-11 (7)", "1400-04-09 (8)", "1401-04-01 (8)", "1402-02-30 (8)", "1402-02-30 (2)")
for ($i = 0; $i -lt $headerRow9.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $headerRow9[$i]
}

# --- Data rows 12-58: shift D:H one period left, new period value lands in H ---
$rowData = @{
    12 = @(6518, 33475, 50630, 124711, 255126)
    13 = @(0, 0, 0, 0, 0)
    14 = @(414799, 480373, 768470, 1298799, 4133142)
    15 = @(150854, 327077, 546798, 680684, 1509807)
    16 = @(11188, 6976, 71616, 693571, 468936)
    17 = @(0, 0, 0, 0, 0)
    18 = @(583359, 847901, 1437514, 2797765, 6367011)
    19 = @(0, 0, 0, 0, 0)
    20 = @(69, 69, 69, 69, 69)
    21 = @(0, 0, 0, 0, 0)
    22 = @(437846, 4587224, 4473707, 4485773, 4440095)
    23 = @(1302, 1302, 1302, 1302, 1302)
    24 = @("-", "-", "-", "-", "-")
    25 = @(0, 0, 0, 0, 0)
    26 = @(439217, 4588595, 4475078, 4487144, 4441466)
    27 = @(1022576, 5436496, 5912592, 7284909, 10808477)
    29 = @(302382, 297137, 214295, 327294, 803973)
    30 = @("-", "-", "-", "-", "-")
    31 = @(0, 0, 0, 0, 0)
    32 = @(39090, 53785, 114998, 147625, 416996)
    33 = @(3070, 3466, 5323, 8819, 11132)
    34 = @(57799, 29267, 42755, 532549, 683417)
    35 = @(0, 0, 0, 0, 0)
    36 = @(0, 0, 0, 0, 0)
    37 = @(402341, 383655, 377371, 1016287, 1915518)
    38 = @(0, 0, 0, 0, 0)
    39 = @("-", "-", "-", "-", "-")
    40 = @(0, 0, 27403, 0, 0)
    41 = @(82595, 102192, 130856, 250921, 406132)
    42 = @(82595, 102192, 158259, 250921, 406132)
    43 = @(484936, 485847, 535630, 1267208, 2321650)
    45 = @(300000, 300000, 4484000, 4484000, 4484000)
    46 = @(0, 0, 0, 0, 0)
    47 = @(0, 0, 0, 0, 0)
    48 = @(0, 0, 0, 0, 0)
    49 = @(0, 0, 0, 0, 0)
    50 = @(24606, 30000, 54006, 97253, 243129)
    51 = @(12895, 12895, 12895, 12895, 12895)
    52 = @("-", "-", "-", "-", "-")
    53 = @(2430, 4202873, 18350, 17827, 17304)
    54 = @("-", "-", "-", "-", "-")
    55 = @(0, 0, 0, 0, 0)
    56 = @(197709, 404881, 807711, 1405726, 3729499)
    57 = @(537640, 4950649, 5376962, 6017701, 8486827)
    58 = @(1022576, 5436496, 5912592, 7284909, 10808477)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
    }
}
